# Commit: Mon, Jul 13, 2020 10:05:25 PM
#
# Change: the table on slide 5 ("B1 and B2" deck) gets a different
# (built-in) table style applied - its tableStyleId switches from the
# custom "Table_0" style ({D1AD5221-5529-4285-8DD7-78F413EDA564}) to the
# built-in style {64757207-97FE-4F56-9239-9357E6D35325}.

$p = $ppt.ActivePresentation

# Locate the shape that holds the table (slide 5, 2nd shape).
$slide = $p.Slides.Item(5)
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $tableShape = $shp
        break
    }
}

$table = $tableShape.Table
$table.ApplyStyle("{64757207-97FE-4F56-9239-9357E6D35325}")
